$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.03147593874317
$ws.Range("D2").Value = 1.040309344572511
$ws.Range("E2").Value = 0.992614727750844
$ws.Range("F2").Value = 1.048795129013439
$ws.Range("I2").Value = 1.035289422311934
$ws.Range("J2").Value = 1.036611456736802
$ws.Range("K2").Value = 1.043092058362471
$ws.Range("L2").Value = 0.9955398523335997
$ws.Range("M2").Value = 1.051553978549312
$ws.Range("N2").Value = 1.038083563906364

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.032417529451702
$ws.Range("D3").Value = 1.041033724713244
$ws.Range("E3").Value = 0.9936372048519299
$ws.Range("F3").Value = 1.049651937296318
$ws.Range("I3").Value = 1.035445818570116
$ws.Range("J3").Value = 1.037194978424624
$ws.Range("K3").Value = 1.043627061932231
$ws.Range("L3").Value = 0.9963617723202687
$ws.Range("M3").Value = 1.052222788321959
$ws.Range("N3").Value = 1.038667914261914

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.033027218632859
$ws.Range("D4").Value = 1.041502644659833
$ws.Range("E4").Value = 0.9942998659930998
$ws.Range("F4").Value = 1.050206876902091
$ws.Range("I4").Value = 1.035545779894282
$ws.Range("J4").Value = 1.03757236278447
$ws.Range("K4").Value = 1.043972772919909
$ws.Range("L4").Value = 0.9968940712668347
$ws.Range("M4").Value = 1.052655436922503
$ws.Range("N4").Value = 1.039045834550838

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.033283630827136
$ws.Range("D5").Value = 1.041699824462591
$ws.Range("E5").Value = 0.994578699834602
$ws.Range("F5").Value = 1.050440298254512
$ws.Range("I5").Value = 1.035587506818886
$ws.Range("J5").Value = 1.03773096794778
$ws.Range("K5").Value = 1.044117995716791
$ws.Range("L5").Value = 0.9971179600053012
$ws.Range("M5").Value = 1.052837293349906
$ws.Range("N5").Value = 1.039204664951672

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.033326689331894
$ws.Range("D6").Value = 1.041732934466344
$ws.Range("E6").Value = 0.994625531979634
$ws.Range("F6").Value = 1.050479497997792
$ws.Range("I6").Value = 1.035594495531634
$ws.Range("J6").Value = 1.037757595655273
$ws.Range("K6").Value = 1.044142372540934
$ws.Range("L6").Value = 0.9971555583673455
$ws.Range("M6").Value = 1.052867826104376
$ws.Range("N6").Value = 1.039231330473563

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.0330306444379
$ws.Range("D7").Value = 1.04150527920763
$ws.Range("E7").Value = 0.994303590798249
$ws.Range("F7").Value = 1.050209995401369
$ws.Range("I7").Value = 1.035546338618465
$ws.Range("J7").Value = 1.037574482261191
$ws.Range("K7").Value = 1.043974713843731
$ws.Range("L7").Value = 0.9968970624462089
$ws.Range("M7").Value = 1.052657867012698
$ws.Range("N7").Value = 1.039047957037459

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.031794066795399
$ws.Range("D8").Value = 1.040554110462704
$ws.Range("E8").Value = 0.9929600610674297
$ws.Range("F8").Value = 1.049084581286144
$ws.Range("I8").Value = 1.035342533062195
$ws.Range("J8").Value = 1.036808700040838
$ws.Range("K8").Value = 1.043272962424485
$ws.Range("L8").Value = 0.9958175282591056
$ws.Range("M8").Value = 1.05178002921196
$ws.Range("N8").Value = 1.038281087318516

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.029618295970935
$ws.Range("D9").Value = 1.038879606089695
$ws.Range("E9").Value = 0.9906006454969559
$ws.Range("F9").Value = 1.047105563218291
$ws.Range("I9").Value = 1.034973948176396
$ws.Range("J9").Value = 1.035457857128428
$ws.Range("K9").Value = 1.042032818747947
$ws.Range("L9").Value = 0.9939188001724441
$ws.Range("M9").Value = 1.050232333786501
$ws.Range("N9").Value = 1.036928326054226

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.028170015073451
$ws.Range("D10").Value = 1.037764414656977
$ws.Range("E10").Value = 0.989033133672735
$ws.Range("F10").Value = 1.045789066966745
$ws.Range("I10").Value = 1.034721901567865
$ws.Range("J10").Value = 1.034556386131972
$ws.Range("K10").Value = 1.041203721008005
$ws.Range("L10").Value = 0.9926553831429383
$ws.Range("M10").Value = 1.049200051122376
$ws.Range("N10").Value = 1.036025574865555

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.027543434580535
$ws.Range("D11").Value = 1.037281813688514
$ws.Range("E11").Value = 0.988355674866747
$ws.Range("F11").Value = 1.045219703837759
$ws.Range("I11").Value = 1.034611269555456
$ws.Range("J11").Value = 1.034165835387548
$ws.Range("K11").Value = 1.040844171842638
$ws.Range("L11").Value = 0.9921088820399291
$ws.Range("M11").Value = 1.048752961586013
$ws.Range("N11").Value = 1.035634469494277

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.027310775888136
$ws.Range("D12").Value = 1.037102598313777
$ws.Range("E12").Value = 0.9881042295826724
$ws.Range("F12").Value = 1.045008321915865
$ws.Range("I12").Value = 1.034569951934899
$ws.Range("J12").Value = 1.03402073713761
$ws.Range("K12").Value = 1.040710538519655
$ws.Range("L12").Value = 0.9919059725120875
$ws.Range("M12").Value = 1.048586878165527
$ws.Range("N12").Value = 1.035489165188182

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.027360678278907
$ws.Range("D13").Value = 1.03714103860132
$ws.Range("E13").Value = 0.9881581567098651
$ws.Range("F13").Value = 1.045053659285418
$ws.Range("I13").Value = 1.034578824831361
$ws.Range("J13").Value = 1.034051862561878
$ws.Range("K13").Value = 1.04073920695879
$ws.Range("L13").Value = 0.9919494934313052
$ws.Range("M13").Value = 1.048622504257855
$ws.Range("N13").Value = 1.035520334814122

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.02752420127892
$ws.Range("D14").Value = 1.037266998781739
$ws.Range("E14").Value = 0.9883348863814464
$ws.Range("F14").Value = 1.045202228793853
$ws.Range("I14").Value = 1.034607858793883
$ws.Range("J14").Value = 1.034153842134059
$ws.Range("K14").Value = 1.040833127311864
$ws.Range("L14").Value = 0.9920921077337197
$ws.Range("M14").Value = 1.048739233366623
$ws.Range("N14").Value = 1.035622459208993

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.027624963982196
$ws.Range("D15").Value = 1.037344612893403
$ws.Range("E15").Value = 0.9884438009545853
$ws.Range("F15").Value = 1.045293781314856
$ws.Range("I15").Value = 1.034625717914334
$ws.Range("J15").Value = 1.034216671124643
$ws.Range("K15").Value = 1.040890984071323
$ws.Range("L15").Value = 0.9921799884222134
$ws.Range("M15").Value = 1.048811152143311
$ws.Range("N15").Value = 1.03568537742395

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.028211610496672
$ws.Range("D16").Value = 1.037796449403098
$ws.Range("E16").Value = 0.9890781214508737
$ws.Range("F16").Value = 1.045826868352114
$ws.Range("I16").Value = 1.034729212413844
$ws.Range("J16").Value = 1.034582301397352
$ws.Range("K16").Value = 1.041227571720735
$ws.Range("L16").Value = 0.9926916645766087
$ws.Range("M16").Value = 1.049229720887997
$ws.Range("N16").Value = 1.036051526933584

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.028579742121388
$ws.Range("D17").Value = 1.038079951609239
$ws.Range("E17").Value = 0.989476357848556
$ws.Range("F17").Value = 1.046161445115143
$ws.Range("I17").Value = 1.034793732112349
$ws.Range("J17").Value = 1.03481159660396
$ws.Range("K17").Value = 1.041438559245176
$ws.Range("L17").Value = 0.9930127773699352
$ws.Range("M17").Value = 1.049492250928328
$ws.Range("N17").Value = 1.036281147765687

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.028794518601379
$ws.Range("D18").Value = 1.038245340984633
$ws.Range("E18").Value = 0.9897087662937556
$ws.Range("F18").Value = 1.046356664406238
$ws.Range("I18").Value = 1.034831221160694
$ws.Range("J18").Value = 1.03494532050241
$ws.Range("K18").Value = 1.041561572105298
$ws.Range("L18").Value = 0.9932001317071769
$ws.Range("M18").Value = 1.049645370111365
$ws.Range("N18").Value = 1.036415061567411

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.028867760544987
$ws.Range("D19").Value = 1.03830173908139
$ws.Range("E19").Value = 0.9897880325774034
$ws.Range("F19").Value = 1.046423240344995
$ws.Range("I19").Value = 1.03484397948282
$ws.Range("J19").Value = 1.034990913456204
$ws.Range("K19").Value = 1.04160350735954
$ws.Range("L19").Value = 0.9932640239640975
$ws.Range("M19").Value = 1.049697578037624
$ws.Range("N19").Value = 1.036460719268429

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.028540239745541
$ws.Range("D20").Value = 1.038049531668015
$ws.Range("E20").Value = 0.9894336180360679
$ws.Range("F20").Value = 1.046125541323875
$ws.Range("I20").Value = 1.034786824670873
$ws.Range("J20").Value = 1.034786997472409
$ws.Range("K20").Value = 1.041415927713062
$ws.Range("L20").Value = 0.9929783193494215
$ws.Range("M20").Value = 1.049464084980402
$ws.Range("N20").Value = 1.036256513700547

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.027476045582566
$ws.Range("D21").Value = 1.037229905419096
$ws.Range("E21").Value = 0.9882828385668249
$ws.Range("F21").Value = 1.045158475869233
$ws.Range("I21").Value = 1.034599315193957
$ws.Range("J21").Value = 1.034123812521306
$ws.Range("K21").Value = 1.040805472329709
$ws.Range("L21").Value = 0.9920501090198102
$ws.Range("M21").Value = 1.048704859942903
$ws.Range("N21").Value = 1.035592386950746

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.026807414712974
$ws.Range("D22").Value = 1.036714829781447
$ws.Range("E22").Value = 0.9875604150241495
$ws.Range("F22").Value = 1.044551051050854
$ws.Range("I22").Value = 1.0344801250061
$ws.Range("J22").Value = 1.033706666884888
$ws.Range("K22").Value = 1.040421188093916
$ws.Range("L22").Value = 0.9914670000341481
$ws.Range("M22").Value = 1.048227421968023
$ws.Range("N22").Value = 1.035174648919677

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.027161823723253
$ws.Range("D23").Value = 1.036987856396726
$ws.Range("E23").Value = 0.9879432794643023
$ws.Range("F23").Value = 1.044873000355845
$ws.Range("I23").Value = 1.034543432609086
$ws.Range("J23").Value = 1.033927819956148
$ws.Range("K23").Value = 1.040624948332219
$ws.Range("L23").Value = 0.991776070289318
$ws.Range("M23").Value = 1.048480528339651
$ws.Range("N23").Value = 1.035396116053664

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.028558089009642
$ws.Range("D24").Value = 1.038063277045726
$ws.Range("E24").Value = 0.9894529299347244
$ws.Range("F24").Value = 1.046141764496898
$ws.Range("I24").Value = 1.034789946291949
$ws.Range("J24").Value = 1.034798112823329
$ws.Range("K24").Value = 1.041426154090858
$ws.Range("L24").Value = 0.9929938892766442
$ws.Range("M24").Value = 1.049476811991349
$ws.Range("N24").Value = 1.036267644836541

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.030180395461615
$ws.Range("D25").Value = 1.039312310033276
$ws.Range("E25").Value = 0.9912096547607049
$ws.Range("F25").Value = 1.047616691152373
$ws.Range("I25").Value = 1.03507035312939
$ws.Range("J25").Value = 1.035807246799326
$ws.Range("K25").Value = 1.042353841882679
$ws.Range("L25").Value = 0.9944092447426414
$ws.Range("M25").Value = 1.050632541206539
$ws.Range("N25").Value = 1.038083563906364

